$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A5 previously used the plain label style; switch it to the bold/italic
# "section label" style already used by the sibling "1st contact" cells
# (A6/A17/A28/A39) by copying that cell's formatting across.
$ws.Range("A6").Copy()
$ws.Range("A5").PasteSpecial(-4122)

# Set the "Time group" label text in column A for the four trial blocks.
$ws.Range("A5").Value = "Time group"
$ws.Range("A16").Value = "Time group"
$ws.Range("A27").Value = "Time group"
$ws.Range("A38").Value = "Time group"

# Update the active selection to match the saved view (A38 single cell).
$ws.Range("A38").Select()
